# Generate Report for Handoff
# Adds two new source files (787154a6-... and 8d2fa0c0-...) to the
# localization status report: one new row per file in the "Overview"
# sheet, and one new row per file in each of the "zh-cn" / "de-de"
# per-locale sheets. Tables / filters / dimensions are resized to match.

$wb = $excel.ActiveWorkbook

$dateFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 6 - 787154a6-d431-4486-ba60-c634e9d6534b.md
$ov.Cells.Item(6,1).Value = "787154a6-d431-4486-ba60-c634e9d6534b.md"
$ov.Hyperlinks.Add($ov.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a1b2c3d4e5f607182930415263748596a7b8c9d/e2e/787154a6-d431-4486-ba60-c634e9d6534b.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\787154a6-d431-4486-ba60-c634e9d6534b.md") | Out-Null
$ov.Cells.Item(6,3).Value = ".md"
$ov.Cells.Item(6,4).Value = "'"
$ov.Cells.Item(6,5).Value = "Ready for handoff"
$ov.Cells.Item(6,6).Value = "Ready for handoff"
$ov.Cells.Item(6,7).NumberFormat = $dateFormat
$ov.Cells.Item(6,7).Value = "2016-08-17 16:41:06"

# Row 7 - 8d2fa0c0-d030-4620-a5f2-4718f521ea22.md
$ov.Cells.Item(7,1).Value = "8d2fa0c0-d030-4620-a5f2-4718f521ea22.md"
$ov.Hyperlinks.Add($ov.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2c3d4e5f607182930415263748596a7b8c9d0e/e2e/8d2fa0c0-d030-4620-a5f2-4718f521ea22.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\8d2fa0c0-d030-4620-a5f2-4718f521ea22.md") | Out-Null
$ov.Cells.Item(7,3).Value = ".md"
$ov.Cells.Item(7,4).Value = "'"
$ov.Cells.Item(7,5).Value = "Ready for handoff"
$ov.Cells.Item(7,6).Value = "Ready for handoff"
$ov.Cells.Item(7,7).NumberFormat = $dateFormat
$ov.Cells.Item(7,7).Value = "2016-08-17 16:41:06"

# Resize the Overview table + autofilter to include the two new rows
$ovTable = $ov.ListObjects.Item(1)
$ovTable.Resize($ov.Range("A1:G7"))

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 6 - 787154a6-d431-4486-ba60-c634e9d6534b.md
$zh.Cells.Item(6,1).Value = "787154a6-d431-4486-ba60-c634e9d6534b.md"
$zh.Hyperlinks.Add($zh.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a1b2c3d4e5f607182930415263748596a7b8c9d/e2e/787154a6-d431-4486-ba60-c634e9d6534b.md", [System.Type]::Missing, [System.Type]::Missing, "787154a6-d431-4486-ba60-c634e9d6534b.md") | Out-Null
$zh.Cells.Item(6,2).Value = ".md"
$zh.Cells.Item(6,3).Value = "Ready for handoff"
$zh.Cells.Item(6,4).Value = "e2e"
$zh.Cells.Item(6,5).Value = "ht"
$zh.Cells.Item(6,6).Value = "'False"
$zh.Cells.Item(6,7).Value = "787154a6-d431-4486-ba60-c634e9d6534b.30f136487c8660c508a7373dcac221699aa83392.zh-cn.xlf"
$zh.Cells.Item(6,8).NumberFormat = $dateFormat
$zh.Cells.Item(6,8).Value = "2016-08-17 16:40:57"
$zh.Cells.Item(6,9).Value = "'"
$zh.Cells.Item(6,10).Value = "'"
$zh.Cells.Item(6,11).NumberFormat = $dateFormat
$zh.Cells.Item(6,11).Value = "0001-01-01 00:00:00"
$zh.Cells.Item(6,12).Value = "'"
$zh.Cells.Item(6,13).Value = "'True"
$zh.Cells.Item(6,14).Value = "'"
$zh.Cells.Item(6,15).Value = "'False"
$zh.Cells.Item(6,16).Value = "'"

# Row 7 - 8d2fa0c0-d030-4620-a5f2-4718f521ea22.md
$zh.Cells.Item(7,1).Value = "8d2fa0c0-d030-4620-a5f2-4718f521ea22.md"
$zh.Hyperlinks.Add($zh.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2c3d4e5f607182930415263748596a7b8c9d0e/e2e/8d2fa0c0-d030-4620-a5f2-4718f521ea22.md", [System.Type]::Missing, [System.Type]::Missing, "8d2fa0c0-d030-4620-a5f2-4718f521ea22.md") | Out-Null
$zh.Cells.Item(7,2).Value = ".md"
$zh.Cells.Item(7,3).Value = "Ready for handoff"
$zh.Cells.Item(7,4).Value = "e2e"
$zh.Cells.Item(7,5).Value = "ht"
$zh.Cells.Item(7,6).Value = "'False"
$zh.Cells.Item(7,7).Value = "8d2fa0c0-d030-4620-a5f2-4718f521ea22.d026c3310eebdfb5a7ee2de9959fd39254edadc8.zh-cn.xlf"
$zh.Cells.Item(7,8).NumberFormat = $dateFormat
$zh.Cells.Item(7,8).Value = "2016-08-17 16:40:57"
$zh.Cells.Item(7,9).Value = "'"
$zh.Cells.Item(7,10).Value = "'"
$zh.Cells.Item(7,11).NumberFormat = $dateFormat
$zh.Cells.Item(7,11).Value = "0001-01-01 00:00:00"
$zh.Cells.Item(7,12).Value = "'"
$zh.Cells.Item(7,13).Value = "'True"
$zh.Cells.Item(7,14).Value = "'"
$zh.Cells.Item(7,15).Value = "'False"
$zh.Cells.Item(7,16).Value = "'"

# Resize the zh-cn table + autofilter to include the two new rows
$zhTable = $zh.ListObjects.Item(1)
$zhTable.Resize($zh.Range("A1:P7"))

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Row 6 - 787154a6-d431-4486-ba60-c634e9d6534b.md
$de.Cells.Item(6,1).Value = "787154a6-d431-4486-ba60-c634e9d6534b.md"
$de.Hyperlinks.Add($de.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a1b2c3d4e5f607182930415263748596a7b8c9d/e2e/787154a6-d431-4486-ba60-c634e9d6534b.md", [System.Type]::Missing, [System.Type]::Missing, "787154a6-d431-4486-ba60-c634e9d6534b.md") | Out-Null
$de.Cells.Item(6,2).Value = ".md"
$de.Cells.Item(6,3).Value = "Ready for handoff"
$de.Cells.Item(6,4).Value = "e2e"
$de.Cells.Item(6,5).Value = "ht"
$de.Cells.Item(6,6).Value = "'False"
$de.Cells.Item(6,7).Value = "787154a6-d431-4486-ba60-c634e9d6534b.30f136487c8660c508a7373dcac221699aa83392.de-de.xlf"
$de.Cells.Item(6,8).NumberFormat = $dateFormat
$de.Cells.Item(6,8).Value = "2016-08-17 16:41:06"
$de.Cells.Item(6,9).Value = "'"
$de.Cells.Item(6,10).Value = "'"
$de.Cells.Item(6,11).NumberFormat = $dateFormat
$de.Cells.Item(6,11).Value = "0001-01-01 00:00:00"
$de.Cells.Item(6,12).Value = "'"
$de.Cells.Item(6,13).Value = "'True"
$de.Cells.Item(6,14).Value = "'"
$de.Cells.Item(6,15).Value = "'False"
$de.Cells.Item(6,16).Value = "'"

# Row 7 - 8d2fa0c0-d030-4620-a5f2-4718f521ea22.md
$de.Cells.Item(7,1).Value = "8d2fa0c0-d030-4620-a5f2-4718f521ea22.md"
$de.Hyperlinks.Add($de.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2c3d4e5f607182930415263748596a7b8c9d0e/e2e/8d2fa0c0-d030-4620-a5f2-4718f521ea22.md", [System.Type]::Missing, [System.Type]::Missing, "8d2fa0c0-d030-4620-a5f2-4718f521ea22.md") | Out-Null
$de.Cells.Item(7,2).Value = ".md"
$de.Cells.Item(7,3).Value = "Ready for handoff"
$de.Cells.Item(7,4).Value = "e2e"
$de.Cells.Item(7,5).Value = "ht"
$de.Cells.Item(7,6).Value = "'False"
$de.Cells.Item(7,7).Value = "8d2fa0c0-d030-4620-a5f2-4718f521ea22.d026c3310eebdfb5a7ee2de9959fd39254edadc8.de-de.xlf"
$de.Cells.Item(7,8).NumberFormat = $dateFormat
$de.Cells.Item(7,8).Value = "2016-08-17 16:41:06"
$de.Cells.Item(7,9).Value = "'"
$de.Cells.Item(7,10).Value = "'"
$de.Cells.Item(7,11).NumberFormat = $dateFormat
$de.Cells.Item(7,11).Value = "0001-01-01 00:00:00"
$de.Cells.Item(7,12).Value = "'"
$de.Cells.Item(7,13).Value = "'True"
$de.Cells.Item(7,14).Value = "'"
$de.Cells.Item(7,15).Value = "'False"
$de.Cells.Item(7,16).Value = "'"

# Resize the de-de table + autofilter to include the two new rows
$deTable = $de.ListObjects.Item(1)
$deTable.Resize($de.Range("A1:P7"))

$ov.Select()
